$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 76092340
$ws.Range("I62").Value = 31257332
$ws.Range("J62").Value = 178572340
$ws.Range("K62").Value = 31257332
$ws.Range("L62").Value = 178572340
$ws.Range("M62").Value = -31256708
$ws.Range("N62").Value = -178573588

# Row 65
$ws.Range("H65").Value = 76092340
$ws.Range("I65").Value = 31257332
$ws.Range("J65").Value = 178572340
$ws.Range("K65").Value = 156286660
$ws.Range("L65").Value = 892861700
$ws.Range("M65").Value = -156283540
$ws.Range("N65").Value = -892867940

# Row 98
$ws.Range("H98").Value = 23717840
$ws.Range("I98").Value = 9525234
$ws.Range("J98").Value = 56833916
$ws.Range("K98").Value = 9525234
$ws.Range("L98").Value = 56833916
$ws.Range("M98").Value = -9523736
$ws.Range("N98").Value = -56836912

# Row 100
$ws.Range("H100").Value = 33337584
$ws.Range("I100").Value = 0
$ws.Range("J100").Value = 33337584
$ws.Range("K100").Value = 0
$ws.Range("L100").Value = 33337584
$ws.Range("M100").ClearContents()
$ws.Range("N100").Value = -33338666

# Row 111
$ws.Range("H111").Value = 13856.2
$ws.Range("I111").Value = 4618.5713
$ws.Range("J111").Value = 35410.668
$ws.Range("K111").Value = 13855.7139
$ws.Range("L111").Value = 106232.004
$ws.Range("M111").Value = -10788.7139
$ws.Range("N111").Value = -112366.004

# Row 122
$ws.Range("H122").Value = 23717840
$ws.Range("I122").Value = 9525234
$ws.Range("J122").Value = 56833916
$ws.Range("K122").Value = 28575702
$ws.Range("L122").Value = 170501748
$ws.Range("M122").Value = -28573252
$ws.Range("N122").Value = -170506648

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 21429812
$ws.Range("I134").Value = 33334250
$ws.Range("J134").Value = 3573157
$ws.Range("K134").Value = 100002750
$ws.Range("L134").Value = 10719471
$ws.Range("M134").Value = -100000215
$ws.Range("N134").Value = -10724541

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2979481.5
$ws.Range("I31").Value = 1603831.1
$ws.Range("J31").Value = 6953583
$ws.Range("K31").Value = 1603831.1
$ws.Range("L31").Value = 6953583
$ws.Range("M31").Value = -1603536.1
$ws.Range("N31").Value = -6954173

# Row 34
$ws.Range("H34").Value = 2979481.5
$ws.Range("I34").Value = 1603831.1
$ws.Range("J34").Value = 6953583
$ws.Range("K34").Value = 1603831.1
$ws.Range("L34").Value = 6953583
$ws.Range("M34").Value = -1603629.1
$ws.Range("N34").Value = -6953987

# Row 99
$ws.Range("H99").Value = 250051250
$ws.Range("I99").Value = 500037500
$ws.Range("J99").Value = 65000
$ws.Range("K99").Value = 500037500
$ws.Range("L99").Value = 65000
$ws.Range("M99").Value = -500036002
$ws.Range("N99").Value = -67996

# Row 122
$ws.Range("H122").Value = 4084.4075
$ws.Range("I122").Value = 5940.8237
$ws.Range("J122").Value = 928.5
$ws.Range("K122").Value = 17822.4711
$ws.Range("L122").Value = 2785.5
$ws.Range("M122").Value = -15372.4711
$ws.Range("N122").Value = -7685.5

# Row 126
$ws.Range("H126").Value = 250051250
$ws.Range("I126").Value = 500037500
$ws.Range("J126").Value = 65000
$ws.Range("K126").Value = 1500112500
$ws.Range("L126").Value = 195000
$ws.Range("M126").Value = -1500110030
$ws.Range("N126").Value = -199940

$ws = $wb.Worksheets.Item("CUL")
# Row 122
$ws.Range("H122").Value = 1611.1333
$ws.Range("I122").Value = 409
$ws.Range("J122").Value = 2663
$ws.Range("K122").Value = 3681
$ws.Range("L122").Value = 23967
$ws.Range("M122").Value = -1231
$ws.Range("N122").Value = -28867

# Row 131
$ws.Range("H131").Value = 50866.55
$ws.Range("I131").Value = 487.36365
$ws.Range("J131").Value = 112441.11
$ws.Range("K131").Value = 1462.09095
$ws.Range("L131").Value = 337323.33
$ws.Range("M131").Value = 3577.90905
$ws.Range("N131").Value = -347403.33

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 8573.817999999999
$ws.Range("I102").Value = 9131.200000000001
$ws.Range("K102").Value = 9131.200000000001
$ws.Range("M102").Value = -7509.200000000001

# Row 122
$ws.Range("H122").Value = 5071340.5
$ws.Range("I122").Value = 34182.473
$ws.Range("J122").Value = 11907484
$ws.Range("K122").Value = 102547.419
$ws.Range("L122").Value = 35722452
$ws.Range("M122").Value = -100097.419
$ws.Range("N122").Value = -35727352

# Row 126
$ws.Range("H126").Value = 14451.75
$ws.Range("I126").Value = 26950
$ws.Range("J126").Value = 1953.5
$ws.Range("K126").Value = 80850
$ws.Range("L126").Value = 5860.5
$ws.Range("M126").Value = -78380
$ws.Range("N126").Value = -10800.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2150
$ws.Range("I7").Value = 2150
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2150
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2038
$ws.Range("N7").ClearContents()

# Row 40
$ws.Range("H40").Value = 111111110
$ws.Range("I40").Value = 111111110
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 111111110
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -111110974
$ws.Range("N40").ClearContents()

# Row 61
$ws.Range("H61").Value = 2500
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 2500
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -2904

# Row 113
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 2500
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 2500
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -6840

# Row 122
$ws.Range("H122").Value = 27724746
$ws.Range("I122").Value = 33911700
$ws.Range("J122").Value = 22225234
$ws.Range("K122").Value = 101735100
$ws.Range("L122").Value = 66675702
$ws.Range("M122").Value = -101732650
$ws.Range("N122").Value = -66680602

# Row 126
$ws.Range("H126").Value = 2150
$ws.Range("I126").Value = 2150
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6450
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3980
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 25147.75
$ws.Range("I107").Value = 33330.332
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 99990.99600000001
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = -98070.99600000001
$ws.Range("N107").Value = -5640

# Row 126
$ws.Range("H126").Value = 71428920
$ws.Range("I126").Value = 71428920
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 214286760
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -214284290
$ws.Range("N126").Value = -214284290
